$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 38288
$ws.Range("D2").Value = 55364299
$ws.Range("C3").Value = 92043
$ws.Range("D3").Value = 134924033
$ws.Range("C4").Value = 31471
$ws.Range("D4").Value = 46606598
$ws.Range("C5").Value = 8816
$ws.Range("D5").Value = 13102856
$ws.Range("C6").Value = 2030
$ws.Range("D6").Value = 3016971
$ws.Range("C12").Value = 41764
$ws.Range("D12").Value = 56659794
$ws.Range("C13").Value = 9794
$ws.Range("D13").Value = 14167091
$ws.Range("C14").Value = 26223
$ws.Range("D14").Value = 38452978
$ws.Range("C15").Value = 8388
$ws.Range("D15").Value = 12447478
$ws.Range("C16").Value = 2174
$ws.Range("D16").Value = 3230935
$ws.Range("C20").Value = 10332
$ws.Range("D20").Value = 13670269
$ws.Range("C21").Value = 13557
$ws.Range("D21").Value = 19566546
$ws.Range("C22").Value = 31937
$ws.Range("D22").Value = 46865995
$ws.Range("C23").Value = 10302
$ws.Range("D23").Value = 15314530
$ws.Range("C24").Value = 2667
$ws.Range("D24").Value = 3965271
$ws.Range("C25").Value = 513
$ws.Range("D25").Value = 763592
$ws.Range("C27").Value = 11801
$ws.Range("D27").Value = 15757678
$ws.Range("C28").Value = 7773
$ws.Range("D28").Value = 11252427
$ws.Range("C29").Value = 22768
$ws.Range("D29").Value = 33418046
$ws.Range("C30").Value = 7885
$ws.Range("D30").Value = 11729391
$ws.Range("C31").Value = 1992
$ws.Range("D31").Value = 2972499
$ws.Range("C34").Value = 8396
$ws.Range("D34").Value = 11091616
$ws.Range("C35").Value = 3300
$ws.Range("D35").Value = 4765653
$ws.Range("C36").Value = 7943
$ws.Range("D36").Value = 11601137
$ws.Range("C37").Value = 3208
$ws.Range("D37").Value = 4754961
$ws.Range("C38").Value = 833
$ws.Range("D38").Value = 1240723
$ws.Range("C41").Value = 2506
$ws.Range("D41").Value = 3389112
$ws.Range("C42").Value = 17486
$ws.Range("D42").Value = 25285547
$ws.Range("C43").Value = 51708
$ws.Range("D43").Value = 75797295
$ws.Range("C44").Value = 19154
$ws.Range("D44").Value = 28450677
$ws.Range("C45").Value = 5663
$ws.Range("D45").Value = 8430805
$ws.Range("C46").Value = 1225
$ws.Range("D46").Value = 1828045
$ws.Range("C50").Value = 16912
$ws.Range("D50").Value = 22493189
$ws.Range("C51").Value = 2087
$ws.Range("D51").Value = 3027046
$ws.Range("C52").Value = 7064
$ws.Range("D52").Value = 10382410
$ws.Range("C53").Value = 2392
$ws.Range("D53").Value = 3572464
$ws.Range("C54").Value = 760
$ws.Range("D54").Value = 1135305
$ws.Range("C57").Value = 7161
$ws.Range("D57").Value = 9847012
$ws.Range("C58").Value = 1063
$ws.Range("D58").Value = 1710006
$ws.Range("C59").Value = 2668
$ws.Range("D59").Value = 4325130
$ws.Range("C60").Value = 1054
$ws.Range("D60").Value = 1713338
$ws.Range("C61").Value = 357
$ws.Range("D61").Value = 582883
$ws.Range("C62").Value = 119
$ws.Range("D62").Value = 200600
$ws.Range("C64").Value = 1557
$ws.Range("D64").Value = 2344263
$ws.Range("C65").Value = 15605
$ws.Range("D65").Value = 22537665
$ws.Range("C66").Value = 45232
$ws.Range("D66").Value = 66184661
$ws.Range("C67").Value = 15854
$ws.Range("D67").Value = 23559052
$ws.Range("C68").Value = 4609
$ws.Range("D68").Value = 6864551
$ws.Range("C69").Value = 939
$ws.Range("D69").Value = 1396668
$ws.Range("C73").Value = 15247
$ws.Range("D73").Value = 20086795
$ws.Range("C74").Value = 53102
$ws.Range("D74").Value = 77281292
$ws.Range("C75").Value = 149568
$ws.Range("D75").Value = 220360313
$ws.Range("C76").Value = 64664
$ws.Range("D76").Value = 96357449
$ws.Range("C77").Value = 20690
$ws.Range("D77").Value = 30915322
$ws.Range("C78").Value = 4921
$ws.Range("D78").Value = 7350401
$ws.Range("C85").Value = 52387
$ws.Range("D85").Value = 71228790
$ws.Range("C86").Value = 4690
$ws.Range("D86").Value = 6797046
$ws.Range("C87").Value = 11745
$ws.Range("D87").Value = 17253620
$ws.Range("C88").Value = 3927
$ws.Range("D88").Value = 5853083
$ws.Range("C89").Value = 1358
$ws.Range("D89").Value = 2029489
$ws.Range("C90").Value = 289
$ws.Range("D90").Value = 431012
$ws.Range("C93").Value = 5483
$ws.Range("D93").Value = 7369903
$ws.Range("C94").Value = 1621
$ws.Range("D94").Value = 2335033
$ws.Range("C95").Value = 5265
$ws.Range("D95").Value = 7755535
$ws.Range("C96").Value = 1960
$ws.Range("D96").Value = 2918976
$ws.Range("C101").Value = 3621
$ws.Range("D101").Value = 4793032
$ws.Range("C102").Value = 678
$ws.Range("D102").Value = 1095649
$ws.Range("C103").Value = 405
$ws.Range("D103").Value = 665597
$ws.Range("C104").Value = 143
$ws.Range("D104").Value = 229160
$ws.Range("C105").Value = 51
$ws.Range("D105").Value = 84000
$ws.Range("C106").Value = 25
$ws.Range("D106").Value = 43500
$ws.Range("C107").Value = 10937
$ws.Range("D107").Value = 15867648
$ws.Range("C108").Value = 29509
$ws.Range("D108").Value = 43349752
$ws.Range("C109").Value = 9867
$ws.Range("D109").Value = 14672889
$ws.Range("C110").Value = 2720
$ws.Range("D110").Value = 4056207
$ws.Range("C111").Value = 496
$ws.Range("D111").Value = 739046
$ws.Range("C114").Value = 9895
$ws.Range("D114").Value = 13069081
$ws.Range("C115").Value = 30899
$ws.Range("D115").Value = 44552386
$ws.Range("C116").Value = 66858
$ws.Range("D116").Value = 97838514
$ws.Range("C117").Value = 21566
$ws.Range("D117").Value = 32051213
$ws.Range("C118").Value = 6116
$ws.Range("D118").Value = 9112021
$ws.Range("C119").Value = 1143
$ws.Range("D119").Value = 1708271
$ws.Range("C124").Value = 26130
$ws.Range("D124").Value = 34887097
$ws.Range("C125").Value = 36543
$ws.Range("D125").Value = 52735403
$ws.Range("C126").Value = 77688
$ws.Range("D126").Value = 113594902
$ws.Range("C127").Value = 24071
$ws.Range("D127").Value = 35726787
$ws.Range("C128").Value = 6443
$ws.Range("D128").Value = 9575238
$ws.Range("C129").Value = 1257
$ws.Range("D129").Value = 1868811
$ws.Range("C133").Value = 32138
$ws.Range("D133").Value = 42665130
$ws.Range("C134").Value = 13456
$ws.Range("D134").Value = 19480832
$ws.Range("C135").Value = 32686
$ws.Range("D135").Value = 48003195
$ws.Range("C136").Value = 11570
$ws.Range("D136").Value = 17190433
$ws.Range("C137").Value = 2989
$ws.Range("D137").Value = 4454741
$ws.Range("C138").Value = 507
$ws.Range("D138").Value = 754490
$ws.Range("C141").Value = 10918
$ws.Range("D141").Value = 14556075
$ws.Range("C142").Value = 35650
$ws.Range("D142").Value = 51489558
$ws.Range("C143").Value = 82314
$ws.Range("D143").Value = 120594719
$ws.Range("C144").Value = 24610
$ws.Range("D144").Value = 36563527
$ws.Range("C145").Value = 6464
$ws.Range("D145").Value = 9645567
$ws.Range("C146").Value = 1462
$ws.Range("D146").Value = 2175230
$ws.Range("C149").Value = 29532
$ws.Range("D149").Value = 39823704
